$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newSource = 'Page 1 (loc: bbox(top_left, xywh): [256, 575, 300, 25]): "Model:20HK2 Rated Capacity:4102mAh, Nominal Energy:64Wh, Limited Charging Voltage:18.04V, Nominal Voltage:15.6V"'

$ws.Range("C4").Value = $newSource
$ws.Range("C5").Value = $newSource
$ws.Range("C8").Value = $newSource
$ws.Range("C9").Value = $newSource

$ws.Range("B13").Value = "頁碼採絕對 1-based。文件主要語言為英文，含少量中文。未提供典型/正常容量與能量，故對應欄位為 null。"
